# This workbook gains one new weekly price record. A new row is inserted
# at row 22 (pushing all subsequent "Haba" records for Femacal de La Calera
# down by one row, so the old row 63 becomes row 64), and the new row 22
# is populated with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 22; this shifts rows 22..63 down
# to 23..64 and extends the sheet dimension to A1:R64 automatically.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new record.
$ws.Cells.Item(22, 1).Value = 3
$ws.Cells.Item(22, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(22, 6).Value = 100112026
$ws.Cells.Item(22, 7).Value = "Haba"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 110
$ws.Cells.Item(22, 11).Value = 9500
$ws.Cells.Item(22, 12).Value = 10000
$ws.Cells.Item(22, 13).Value = 9773
$ws.Cells.Item(22, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 16).Value = 391
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
